# Generate Report for Handoff
# Update the localization-status report so the file
# "d6875cf9-4019-4dc6-b293-0d94272724f5.md" reflects that it is now
# "Ready for handoff" (instead of "Handed back: in sync with en-US"),
# with refreshed handoff timestamps and an error detail message about
# the handback file being stale.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19335bd6e5fb53c334e800bc0b4f8b07ac8bf5c4/e2e/d6875cf9-4019-4dc6-b293-0d94272724f5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e7505039b4b6f21eb9340cd2881d9a403b056679/e2e/d6875cf9-4019-4dc6-b293-0d94272724f5.md."

# --- Overview sheet: row 3 is the d6875cf9-... file ---
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-19 18:56:12"

# --- zh-cn sheet: row 3 is the d6875cf9-... file ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-19 18:56:01"
$zhcn.Range("P3").Value = $errorDetail
# widen the "Error Detail" column to fit the new long message (saved width="40")
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 is the d6875cf9-... file ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-19 18:56:12"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.17

Write-Host "Updated handoff status for d6875cf9-4019-4dc6-b293-0d94272724f5.md"
